$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Swap values in columns A, Q, R between row 2 and row 3
$colsToSwap = @("A", "Q", "R")

foreach ($col in $colsToSwap) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $val2 = $cell2.Value()
    $val3 = $cell3.Value()

    $cell2.Value = $val3
    $cell3.Value = $val2
}
